# Apply the commit's changes to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet (shown in workbook.xml sheet name) from
#    "RGossF-HW30.xpc" to "RGossF".
$ws.Name = "RGossF"

# 2. Tiny floating point tweaks on row 13 (rounding refresh from the
#    Gaussian Quadrature computation).
$ws.Range("D13").Value = 0.9955323564408874
$ws.Range("J13").Value = 0.9955323564408874
$ws.Range("K13").Value = 0.9955324820543541
$ws.Range("L13").Value = 0.9952503594389865
$ws.Range("M13").Value = 0.9954648634716651

# 3. Add a new row (14th scheme / "HexGrid-60degTilt5degRes") of averaged
#    intensities as row 16.
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9978792165590392
$ws.Range("D16").Value = 0.9568473323611687
$ws.Range("E16").Value = 1.010399892071232
$ws.Range("F16").Value = 0.9978792165590392
$ws.Range("G16").Value = 0.9686704768203884
$ws.Range("H16").Value = 1.038396101144211
$ws.Range("I16").Value = 1.008068444588059
$ws.Range("J16").Value = 0.9568473323611687
$ws.Range("K16").Value = 0.9836236122162003
$ws.Range("L16").Value = 0.9907514143876198
$ws.Range("M16").Value = 0.9967102439240163
